$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controller - testing")

# Row 8: new row for "admins_controller" - testing just started, nothing filled in yet.
$ws.Range("A7").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "admins_controller"

# Row 6 ("user controller"): user testing is done - same destroy bug, plus update bug found.
$ws.Range("D6").Value = "no"
$ws.Range("E6").Value = "destroy method is missing, update method does not work"

# Row 7 ("patients controller"): patients testing started, same findings so far.
$ws.Range("D7").Value = "no"
$ws.Range("E7").Value = "destroy method is missing, update method does not work"
$ws.Rows.Item(7).RowHeight = 45

# Match the printed page orientation recorded for this sheet.
$ws.PageSetup.Orientation = 1

$ws.Range("B8").Select()
